$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shift the existing telefono/direccion/monto columns one column to the
# right (C->D, D->E, E->F) for every header/data row, to make room for the
# new "tipo" column at C. Work right-to-left (E->F first, then D->E, then
# C->D) so that a column's old content is copied out before it gets
# overwritten. Using Copy(destination) (instead of reading/writing .Value)
# also carries the original cell style/number format along automatically.
# ---------------------------------------------------------------------------

$rows = 1, 2, 3, 4, 5

foreach ($r in $rows) {
    $ws.Cells.Item($r, 5).Copy($ws.Cells.Item($r, 6))   # E -> F (monto)
    $ws.Cells.Item($r, 4).Copy($ws.Cells.Item($r, 5))   # D -> E (direccion)
    $ws.Cells.Item($r, 3).Copy($ws.Cells.Item($r, 4))   # C -> D (telefono)
}

# ---------------------------------------------------------------------------
# Header row (row 1): rename "cliente" to "name", and add the "tipo" and
# "anotaciones" headers.
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "tipo"
$ws.Range("C1").ClearFormats()
$ws.Range("G1").Value = "anotaciones"
$ws.Range("G1").ClearFormats()

# ---------------------------------------------------------------------------
# New "tipo" column values - every cliente is of type CUENTA.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "CUENTA"
$ws.Range("C2").ClearFormats()
$ws.Range("C3").Value = "CUENTA"
$ws.Range("C4").Value = "CUENTA"
$ws.Range("C5").Value = "CUENTA"

# ---------------------------------------------------------------------------
# New "anotaciones" column values.
# ---------------------------------------------------------------------------
$ws.Range("G2").Value = "el mascapito"
$ws.Range("G2").ClearFormats()
$ws.Range("G3").Value = "el mascapoto"
$ws.Range("G4").Value = "el capo de la mafia"
$ws.Range("G5").Value = "chupapijas"

# ---------------------------------------------------------------------------
# Touch the next free cell below/right of the table, growing the used range
# down to row 11 / column H (carrying over the same blank placeholder style
# already used by E8/E9/G10), and select it - mirrors the author clicking
# past the last filled-in cell after finishing their edits.
# ---------------------------------------------------------------------------
$ws.Range("G10").Copy($ws.Range("H11"))
$ws.Rows.Item(11).RowHeight = 15.75
$null = $ws.Range("H11").Select()
